# Generate Report for Handoff
#
# Refreshes the localization-status report after a new handoff run:
#  - Overview sheet: "Latest HO Xliff Generate Date" (col G) moves forward
#    16s for every file that was part of this handoff batch.
#  - zh-cn sheet: "Latest Handoff Datetime" (col H) moves forward 16s for
#    the same rows, and "Priority" (col E) is now stamped "ht".
#  - de-de sheet: "Priority" (col E) is now stamped "ht" for the same rows
#    (its handoff datetime was not part of this batch).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$rows = @(7, 8, 10, 12, 13, 14)

foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-20 22:19:51"

    $wsZhCn.Range("H$r").Value = "2016-08-20 22:19:47"
    $wsZhCn.Range("E$r").Value = "ht"

    $wsDeDe.Range("E$r").Value = "ht"
}
